$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update passenger details in row 2 (values changed; headers in row 1 stay the same)
$ws.Range("A2").Value = "Vamsi"
$ws.Range("B2").Value = "Yellamraju"
$ws.Range("C2").Value = "MALE"
$ws.Range("D2").Value = "7406683580"
$ws.Range("E2").Value = "yvamsipanda9@gmail.com"

# Update the saved selection/active cell to B3
$ws.Range("B3").Select()
